$wb = $excel.ActiveWorkbook
$mcc = $wb.Worksheets.Item("MCC")
$median = $wb.Worksheets.Item("acc_median")
$mcc.Copy($null, $median)
$new = $wb.Worksheets.Item("MCC (2)")
$new.Name = "acc_std"

$new.Cells.Item(3,3).Value = 0.025014943966047601
$new.Cells.Item(3,4).Value = 0.098075587692774996
$new.Cells.Item(3,5).Value = 0.23903381660265699
$new.Cells.Item(3,6).Value = 0.083178879031959099
$new.Cells.Item(3,7).Value = 0.098722514466403705
$new.Cells.Item(3,8).Value = 0.083499914365174704
$new.Cells.Item(3,9).Value = 0.086402204739379901
$new.Cells.Item(3,10).Value = 0.0643209353935211

$new.Cells.Item(4,3).Value = 0.14262807441350001
$new.Cells.Item(4,4).Value = 0.14953852658674699
$new.Cells.Item(4,5).Value = 0.108116754983799
$new.Cells.Item(4,6).Value = 0.074873871776447207
$new.Cells.Item(4,7).Value = 0.16148958980900899
$new.Cells.Item(4,8).Value = 0.16711238246454599
$new.Cells.Item(4,9).Value = 0.12097998767271501
$new.Cells.Item(4,10).Value = 0.173934251701958

$new.Cells.Item(5,3).Value = 0.12138285804827
$new.Cells.Item(5,4).Value = 0.072212816483470901
$new.Cells.Item(5,5).Value = 0.053215203622950701
$new.Cells.Item(5,6).Value = 0.068063570900894405
$new.Cells.Item(5,7).Value = 0.034338956657989403
$new.Cells.Item(5,8).Value = 0.076335634980169501
$new.Cells.Item(5,9).Value = 0.073942530427230896
$new.Cells.Item(5,10).Value = 0.087806973294354995

$new.Cells.Item(6,3).Value = 0.020482660344966999
$new.Cells.Item(6,4).Value = 0.0868168204010802
$new.Cells.Item(6,5).Value = 0.059903507806109403
$new.Cells.Item(6,6).Value = 0.090830884558982194
$new.Cells.Item(6,7).Value = 0.096440908999504693
$new.Cells.Item(6,8).Value = 0.10644065072987199
$new.Cells.Item(6,9).Value = 0.10725612299000099
$new.Cells.Item(6,10).Value = 0.079679958379958601

$new.Cells.Item(7,3).Value = 0.14949015368152899
$new.Cells.Item(7,4).Value = 0.16044457664066999
$new.Cells.Item(7,5).Value = 0.118621109373117
$new.Cells.Item(7,6).Value = 0.089081672559562897
$new.Cells.Item(7,7).Value = 0.12421034710588599
$new.Cells.Item(7,8).Value = 0.13233818967984701
$new.Cells.Item(7,9).Value = 0.111005613663358
$new.Cells.Item(7,10).Value = 0.148471831245376

$new.Cells.Item(8,3).Value = 0.115946368500179
$new.Cells.Item(8,4).Value = 0.032601938270222203
$new.Cells.Item(8,5).Value = 0.060776460516971803
$new.Cells.Item(8,6).Value = 0.081250182765519402
$new.Cells.Item(8,7).Value = 0.136432799899491
$new.Cells.Item(8,8).Value = 0.12518169974641499
$new.Cells.Item(8,9).Value = 0.030453253420055601
$new.Cells.Item(8,10).Value = 0.033138569484302702

$new.Range("C3").Select()

foreach ($s in $wb.Worksheets) {
    Write-Output ($s.Name + " " + $s.Index)
}
